$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New leading "id" column: overwrite A2:A16 (previously airline code "VX")
#    with the new numeric value 618 used by the new pax feed format.
# ---------------------------------------------------------------------------
$ws.Range("A2:A16").Value = 618

# ---------------------------------------------------------------------------
# 2. Apply the built-in "Hyperlink" cell style to column O (the link/ref
#    column). Adding-then-removing a real hyperlink is the cleanest way to
#    pull in Excel's native Hyperlink style (theme-coloured, underlined
#    font) without leaving a live hyperlink behind.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("O2"), "", "O2", "", "1")
$ws.Hyperlinks.Delete()

$ws.Range("O2").Copy()
$ws.Range("O3:O16").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. New blank template rows (17-30) below the data, pre-formatted the same
#    way as the data rows: column K keeps the short-date format, column O
#    keeps the Hyperlink style, ready for more rows to be appended later.
# ---------------------------------------------------------------------------
$ws.Range("K2").Copy()
$ws.Range("K17:K29").PasteSpecial(-4122)

$ws.Range("O2").Copy()
$ws.Range("O17:O30").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4. Column widths: best-fit G, K and P to their (now wider) contents.
# ---------------------------------------------------------------------------
$ws.Columns("G:G").AutoFit()
$ws.Columns("K:K").AutoFit()
$ws.Columns("P:P").AutoFit()

# ---------------------------------------------------------------------------
# 5. View state: scroll back to the left edge, roughly to row 10, and select
#    the pax-number column of the existing data.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("A3:A16").Select()
